# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
# Both sheets share the same row layout / values, so the same set of
# row -> new-value updates applies to each of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 463
    11 = 578
    12 = 30
    13 = 303
    18 = 10
    22 = 950
    23 = 1405
    24 = 301
    25 = 331
    27 = 76
    32 = 254
    34 = 1625
    40 = 3702
    42 = 204
    43 = 914
    46 = 70
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
